$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.245.66"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "2.587.47"
$ws.Range("E3").Value = "  -2.32%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  -3.74%  "

$ws.Range("D9").Value = "2.596.50"
$ws.Range("E9").Value = "  -2.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.39%  "

$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("E12").Value = "  -0.66%  "

$ws.Range("E13").Value = "  +1.63%  "

$ws.Range("D14").Value = "3.043.39"
$ws.Range("E14").Value = "  -2.68%  "

$ws.Range("D15").Value = "60.229.49"
$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").Value = "2.594.84"
$ws.Range("E18").Value = "  -2.94%  "

$ws.Range("E19").Value = "  -1.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "352.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("E22").Value = "  -1.20%  "

$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("E26").Value = "  -0.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("E28").Value = "  -3.92%  "

$ws.Range("E29").Value = "  -2.86%  "

$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("E31").Value = "  -1.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.76%  "

$ws.Range("E33").Value = "  -1.24%  "

$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("E36").Value = "  -3.05%  "

$ws.Range("E37").Value = "  +3.56%  "

$ws.Range("E38").Value = "  -2.78%  "

$ws.Range("E39").Value = "  +1.80%  "

$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.836"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "294.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.05%  "

$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("E44").Value = "  -4.79%  "

$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0551"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.16%  "

$ws.Range("E48").Value = "  -4.12%  "

$ws.Range("E49").Value = "  -1.80%  "

$ws.Range("E50").Value = "  +0.28%  "

$ws.Range("D51").Value = "1.986.78"
$ws.Range("E51").Value = "  -2.46%  "
